# issue #5: stock data output to json file
# Insert a new "property_category" column (with value "stock") into the
# 股票 (Stock) worksheet, between the existing "total" and "date" columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(5)

# Shift columns H:J (date, legislator_name, legislator_id) one column to
# the right to make room for the new property_category column at H.
$ws.Columns.Item(8).Insert()

# New header for the inserted column.
$ws.Range("H1").Value = "property_category"

# Populate the new column for each data row with the stock category value.
$ws.Range("H2").Value = "stock"
$ws.Range("H3").Value = "stock"
$ws.Range("H4").Value = "stock"
